## edit.ps1
## Applies the commit "Add files via upload" to ISYS3001.docx:
##   - Replaces the lone "..." (ellipsis) paragraph near the end of the
##     document with two runs of text: the new sentence, and a separate
##     run containing the trailing full stop.
##   - Removes the now-superfluous empty trailing paragraph that used to
##     follow the ellipsis paragraph, so the new text paragraph is
##     immediately followed by the section properties.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Locate the paragraph that currently contains only the ellipsis
#    character ("…") and rewrite its contents as two runs:
#       "The main function ... for users to use"   +   "."
#    We use Range.InsertXML (the documented Word COM mechanism for
#    inserting/overwriting rich content as OOXML) so that the resulting
#    paragraph contains two distinct <w:r> elements rather than being
#    merged into a single run, exactly like the target markup.
# ---------------------------------------------------------------------

$ellipsis = [string][char]8230   # "…"

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    $text = [string]$candidate.Range.Text
    if ($text.Contains($ellipsis)) {
        $target = $candidate
    }
}

if ($target -ne $null) {
    $targetRange = $target.Range

    $sentence = "The main function of version control is to track file " + `
                "changes, which makes it very simple and convenient for " + `
                "users to use"

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
               '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
                   '<pkg:xmlData>' + `
                       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
                           '<w:body>' + `
                               '<w:p>' + `
                                   '<w:r><w:t>' + $sentence + '</w:t></w:r>' + `
                                   '<w:r><w:t>.</w:t></w:r>' + `
                               '</w:p>' + `
                           '</w:body>' + `
                       '</w:document>' + `
                   '</pkg:xmlData>' + `
               '</pkg:part>' + `
           '</pkg:package>'

    $targetRange.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 2. Remove the empty trailing paragraph that used to sit between the
#    (former ellipsis) paragraph and the section properties. In the
#    original document this was an entirely empty <w:p/>. Deleting the
#    paragraph mark that separates it from the previous paragraph joins
#    it away, leaving the text paragraph immediately before <w:sectPr>.
# ---------------------------------------------------------------------

$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastText = [string]$lastPara.Range.Text

# Paragraph.Range.Text always includes the trailing paragraph mark, so an
# empty paragraph reports a length of 1 (just the mark character).
if ($lastText.Length -le 1) {
    $docEnd = $d.Content.End
    $markRange = $d.Range($docEnd - 2, $docEnd - 1)
    $markRange.Delete()
}
